$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45899
$ws.Range("B2").Value = 91.56
$ws.Range("C2").Value = 86.64
$ws.Range("D2").Value = 87.56
$ws.Range("E2").Value = 83.48
$ws.Range("F2").Value = 80.88
$ws.Range("G2").Value = 80.64
$ws.Range("H2").Value = 85.06
$ws.Range("I2").Value = 88.29000000000001
$ws.Range("J2").Value = 77
$ws.Range("K2").Value = 32
$ws.Range("L2").Value = 9.5
$ws.Range("M2").Value = 1.99
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = -0.01
$ws.Range("T2").Value = 8.16
$ws.Range("U2").Value = 57.97
$ws.Range("V2").Value = 89.93000000000001
$ws.Range("W2").Value = 107.12
$ws.Range("X2").Value = 101.12
$ws.Range("Y2").Value = 95
$ws.Range("Z2").Value = 52.66
$ws.Range("AB2").Value = 98.29000000000001
$ws.Range("AD2").Value = 98.52
$ws.Range("AF2").Value = 98.06
